$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3:H3").Value = 1
$ws.Range("I3:J3").Value = 1

$ws.Range("I6").Value = 2

$ws.Range("G10:J10").Value = 2

$ws.Range("J16").Value = 2

$ws.Range("E19:J19").Value = 2

$ws.Range("J27").Value = 2

$ws.Range("I29").Value = 2

$ws.Range("J30").Value = 2
